# Sunbeam.xlsx -- "Add files via upload": the Data sheet's sample gets
# extended from 16 rows to 51 rows (the original 16 physician/monitor
# systolic-BP pairs, re-appended/extended), and the workbook is resaved
# with the Data sheet as the active/selected sheet (scrolled down near the
# bottom of the new data, selection spanning the whole data body).

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")

# Additional Physician / Monitor reading pairs appended as rows 17-51.
$newRows = @(
    @(120, 138),
    @(111, 123),
    @(123, 119),
    @(114, 122),
    @(121, 116),
    @(132, 131),
    @(139, 116),
    @(141, 123),
    @(120, 138),
    @(99, 123),
    @(128, 119),
    @(118, 122),
    @(116, 116),
    @(123, 119),
    @(114, 122),
    @(111, 123),
    @(123, 119),
    @(114, 122),
    @(121, 116),
    @(132, 131),
    @(139, 116),
    @(141, 123),
    @(120, 138),
    @(99, 123),
    @(128, 119),
    @(112, 126),
    @(109, 108),
    @(139, 116),
    @(141, 123),
    @(120, 138),
    @(99, 123),
    @(128, 119),
    @(118, 122),
    @(116, 116),
    @(120, 118)
)

$row = 17
foreach ($pair in $newRows) {
    $wsData.Cells.Item($row, 1).Value = $pair[0]
    $wsData.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Make "Data" the active/selected sheet (the "Description" sheet was the
# active tab before the upload; afterwards "Data" is), scrolled so the
# newly-added tail of the table is in view, with the full data body
# selected.
$wsData.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$wsData.Range("A2:B51").Select()
